$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.494.58"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "3.494.65"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'605.33"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").Value = "'151.01"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("D7").Value = "3.492.53"
$ws.Range("E7").Value = "  -1.41%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D11").Value = "'7.59"
$ws.Range("E11").Value = "  +6.93%  "

$ws.Range("D12").Value = "'0.432"
$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").Value = "'0.0000216"
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("D14").Value = "'32.07"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").Value = "4.080.38"
$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.570.14"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.479.91"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").Value = "'15.43"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "'9.91"
$ws.Range("E21").Value = "  +2.35%  "

$ws.Range("D22").Value = "'446.85"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("D23").Value = "'0.627"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'79.31"
$ws.Range("E24").Value = "  +2.19%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.629.78"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'0.0000127"
$ws.Range("E27").Value = "  -4.47%  "

$ws.Range("D28").Value = "'8.68"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'9.95"
$ws.Range("E29").Value = "  -3.40%  "

$ws.Range("E30").Value = "  -1.49%  "

$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").Value = "'25.62"
$ws.Range("E34").Value = "  -1.45%  "

$ws.Range("D35").Value = "'6.13"
$ws.Range("E35").Value = "  -1.58%  "

$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").Value = "3.485.45"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").Value = "'8.01"
$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +5.01%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").Value = "'177.06"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").Value = "'0.0898"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").Value = "'5.43"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").Value = "'30.28"
$ws.Range("E46").Value = "  +5.91%  "

$ws.Range("D47").Value = "'46.48"
$ws.Range("E47").Value = "  +2.40%  "

$ws.Range("D48").Value = "'1.29"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").Value = "'2.54"
$ws.Range("E49").Value = "  -5.67%  "

$ws.Range("D50").Value = "'7.61"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.989"
$ws.Range("E51").Value = "  -1.16%  "
